$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.654.95"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.281.44"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.79"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.09"
$ws.Range("E6").Value = "  -2.04%  "

$ws.Range("E7").Value = "  -2.42%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -3.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.21"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.33"
$ws.Range("E12").Value = "  +3.51%  "

$ws.Range("E13").Value = "  +0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.635.62"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.280.76"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.602.27"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("E21").Value = "  -2.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.14"
$ws.Range("E22").Value = "  -1.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.66"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("E24").Value = "  -1.82%  "

$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.41"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.91"
$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.53"
$ws.Range("E36").Value = "  -3.20%  "

$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0691"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("E40").Value = "  -2.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.007.48"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0278"
$ws.Range("E44").Value = "  -2.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.04"
$ws.Range("E45").Value = "  +2.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.99"
$ws.Range("E46").Value = "  -3.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -6.05%  "

$ws.Range("E48").Value = "  -2.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  +4.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.60"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.502.74"
$ws.Range("E51").Value = "  -0.98%  "
